$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): values are stored as Text (e.g. "241.63"), not Number, in this sheet.
# Force each target cell to Text format before assigning so the numeric-looking string is not
# auto-converted to a Number by Excel (NumberFormat must be set per-cell; a multi-area Range
# union only applies it to the first area).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "241.63"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.38"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.169"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05529"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.364"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.319"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8042"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9550"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1378"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07299"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03022"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03065"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09303"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.574"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001650"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04699"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005758"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006439"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004982"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001042"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.768"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.102"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3240"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1289"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1950"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0003104"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03837"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006908"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003061"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008169"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005939"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.6834"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1132"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01012"

# --- Plain text cell updates (coin names, links, rank labels)
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("B27").Value = "AAXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("E27").Value = "26AAXTokenAAB"
$ws.Range("B28").Value = "UpBots"
$ws.Range("C28").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("E28").Value = "27UpBotsUBXT"
$ws.Range("E47").Value = "46ACDXExchangeACXTWorstin24h"
